$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CityResaleNum")

$row = 90

# Columns A (Date) and D (Week) must stay as literal text (e.g. "2025-02-26" and "08")
# rather than being auto-converted by Excel into a date serial / number. Force the
# cell to Text format while assigning, then restore the default "Normal" style so
# the new row keeps the same (unstyled) look as the other data rows.
$cellA = $ws.Cells.Item($row, 1)
$cellD = $ws.Cells.Item($row, 4)
$cellA.NumberFormat = "@"
$cellD.NumberFormat = "@"

$cellA.Value = "2025-02-26"
$ws.Cells.Item($row, 2).Value = "08:53:25"
$ws.Cells.Item($row, 3).Value = "Wednesday"
$cellD.Value = "08"

$cellA.Style = "Normal"
$cellD.Style = "Normal"

$ws.Cells.Item($row, 5).Value = 130661
$ws.Cells.Item($row, 6).Value = 141907
$ws.Cells.Item($row, 7).Value = 172454
$ws.Cells.Item($row, 8).Value = 159560
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 146404
$ws.Cells.Item($row, 11).Value = -1
$ws.Cells.Item($row, 12).Value = -1
$ws.Cells.Item($row, 13).Value = 193759
$ws.Cells.Item($row, 14).Value = 115384
$ws.Cells.Item($row, 15).Value = 46658
$ws.Cells.Item($row, 16).Value = 29467
$ws.Cells.Item($row, 17).Value = 69050
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 48481
$ws.Cells.Item($row, 20).Value = -1
